$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new entire column at G, shifting the existing "Brodmann Area" table
# (columns G:K) and the second table's columns (G, I) one place to the right.
$ws.Columns("G:G").Insert()

# New column header for the first table (row 2) plus the per-cluster
# percentage-of-cluster-from-Yeo-masks values (rows 3-16).
$ws.Range("G2").Value = "BA (from Yale BioImage Suite package)"
$ws.Range("G3").Value = 6
$ws.Range("G4").Value = 10
$ws.Range("G5").Value = 10
$ws.Range("G6").Value = 4
$ws.Range("G7").Value = 9
$ws.Range("G8").Value = 6
$ws.Range("G9").Value = 23
$ws.Range("G10").Value = 6
$ws.Range("G11").Value = 6
$ws.Range("G12").Value = 8
$ws.Range("G13").Value = 7
$ws.Range("G14").Value = 11
$ws.Range("G15").Value = 8
$ws.Range("G16").Value = 9

# The old Brodmann-Area "N/A" note for cluster 11 (now shifted into H13) gets
# annotated with more detail.
$ws.Range("H13").Value = "N/A (but 32 from Joe?)"

# The column insert also nudged the second table's (now-empty) G column;
# clear the leftover placeholder cell at G19 so no stray cell remains there.
$ws.Range("G19").Clear()

# Leave the selection where the author ended up after entering the data.
$ws.Range("G17").Select()
